# Add flow-chart branch labels / wrap-up text to the "系統架構圖" (system
# architecture diagram) explanation paragraph.
#
# Before:
#   使用者提問問題>Bot分析問題>如查詢道具類問題則回覆合適的答案>如其他問題則透過爬蟲
#   搜尋與篩選後，回復合適的答案>若無法分析出則再次詢問/告知使用者無法分析
#
# After:
#   使用者提問問題>Bot分析問題>(分支一)如查詢道具類問題則回覆合適的答案>如其他問題則透過爬蟲
#   搜尋與篩選後，回復合適的答案。(分支二)若無法分析出則再次詢問/告知使用者無法分析>回到提問問題。
#
# Applied as three text edits. They are executed right-to-left (last edit in
# the paragraph first) purely so each `Find` still locates its target text
# exactly as it reads in the original document.

$d = $word.ActiveDocument

# 3) Append "回到提問問題。" (preceded by its own ">") right after the final
#    "告知使用者無法分析" in the paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("告知使用者無法分析", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if (-not $found) { throw "could not find '告知使用者無法分析'" }
$rng.Collapse(0)
$rng.InsertAfter(">回到提問問題。")

# 2) Turn the ">" that used to separate "...回復合適的答案" from
#    "若無法分析出..." into "。(分支二)" — closing the first sentence and
#    labelling the second branch.
$rng = $d.Content
$found = $rng.Find.Execute("回復合適的答案>若無法分析出", $true, $false, $false, $false, $false,
                            $true, 1, $false, "回復合適的答案。(分支二)若無法分析出", 2)
if (-not $found) { throw "could not find '回復合適的答案>若無法分析出'" }

# 1) Label the first branch: insert "(分支一)" right before
#    "如查詢道具類問題則回覆合適的答案".
$rng = $d.Content
$found = $rng.Find.Execute("如查詢道具類問題則回覆合適的答案", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if (-not $found) { throw "could not find '如查詢道具類問題則回覆合適的答案'" }
$rng.Collapse(1)
$rng.InsertBefore("(分支一)")
